$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Victoria)
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("K2").Value = 0

# Row 4 (Poblacion)
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = 0

# Row 8 (Bayanan Main)
$ws.Range("B8").Value = 0

# Row 13 (Buli)
$ws.Range("B13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("K13").Value = 1

# Row 15 (B.Silang)
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("K15").Value = 0

# Row 17 (Total)
$ws.Range("B17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0
